# Generate Report for Handback
#
# The localization "handback" cycle completed for both the zh-cn and de-de
# targets: the handed-back files are now in sync with en-US, so the
# per-language status flips from "Ready for handoff" to
# "Handed back: in sync with en-US", the "Latest Handback DateTime" is
# refreshed to the moment of this handback, and the stale
# "handback file is not the latest" error is cleared out.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns (E, F) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText

# --- zh-cn sheet: Status / Latest Handback DateTime / Error Detail ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $statusText
$zhcn.Range("K2").Value = "2016-09-03 04:54:27"
$zhcn.Range("P2").Value = ""

# --- de-de sheet: Status / Latest Handback DateTime / Error Detail ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $statusText
$dede.Range("K2").Value = "2016-09-03 04:54:34"
$dede.Range("P2").Value = ""

# --- Column widths: the Status / Error Detail columns got wider now that
#     the content they hold changed shape (longer status text, emptied
#     error column) ---
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

$zhcn.Columns.Item(3).ColumnWidth = 29.9777047293527
$zhcn.Columns.Item(16).ColumnWidth = 13.7470528738839

$dede.Columns.Item(3).ColumnWidth = 29.9777047293527
$dede.Columns.Item(16).ColumnWidth = 13.7470528738839
